$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("input")

$ws.Range("G1").Value = "distribution"
$ws.Range("C1").Value = "unit"
$ws.Range("B1").Value = "description"
$ws.Range("A1").Value = "parameter"

$ws.Range("A2").Select()
